# Apply GenX signal refresh: Active Signals, Summary Dashboard, Signal History
$wb = $excel.ActiveWorkbook

# BUY => green fill (C6EFCE), SELL => red fill (FFC7CE) -- same palette as before,
# expressed as packed BGR integers for Interior.Color (COM packs 0xBBGGRR).
$greenColor = 198 + (239*256) + (206*65536)   # C6EFCE
$redColor   = 255 + (199*256) + (206*65536)   # FFC7CE

# ===================== Sheet: Active Signals =====================
$ws1 = $wb.Worksheets.Item("Active Signals")

# Row 2: AUDUSD SELL
$ws1.Cells.Item(2,1).Value = "2025-07-28 19:46"
$ws1.Cells.Item(2,2).Value = "AUDUSD"
$ws1.Cells.Item(2,3).Value = "SELL"
$ws1.Cells.Item(2,4).Value = 0.65914
$ws1.Cells.Item(2,5).Value = 0.66237
$ws1.Cells.Item(2,6).Value = 0.6515
$ws1.Cells.Item(2,7).Value = 0.1
$ws1.Cells.Item(2,8).NumberFormat = "@"
$ws1.Cells.Item(2,8).Value = "74.0%"
$ws1.Cells.Item(2,9).Value = 2.37
$ws1.Cells.Item(2,10).Value = "Active"
$ws1.Cells.Item(2,3).Interior.Color = $redColor

# Row 3: USDCAD SELL
$ws1.Cells.Item(3,1).Value = "2025-07-28 19:34"
$ws1.Cells.Item(3,2).Value = "USDCAD"
$ws1.Cells.Item(3,3).Value = "SELL"
$ws1.Cells.Item(3,4).Value = 1.36221
$ws1.Cells.Item(3,5).Value = 1.36531
$ws1.Cells.Item(3,6).Value = 1.3559
$ws1.Cells.Item(3,7).Value = 0.07
$ws1.Cells.Item(3,8).NumberFormat = "@"
$ws1.Cells.Item(3,8).Value = "86.0%"
$ws1.Cells.Item(3,9).Value = 2.03
$ws1.Cells.Item(3,10).Value = "Active"
$ws1.Cells.Item(3,3).Interior.Color = $redColor

# Row 4: NZDUSD SELL
$ws1.Cells.Item(4,1).Value = "2025-07-28 19:54"
$ws1.Cells.Item(4,2).Value = "NZDUSD"
$ws1.Cells.Item(4,3).Value = "SELL"
$ws1.Cells.Item(4,4).Value = 0.58913
$ws1.Cells.Item(4,5).Value = 0.59282
$ws1.Cells.Item(4,6).Value = 0.57921
$ws1.Cells.Item(4,7).Value = 0.08
$ws1.Cells.Item(4,8).NumberFormat = "@"
$ws1.Cells.Item(4,8).Value = "89.0%"
$ws1.Cells.Item(4,9).Value = 2.69
$ws1.Cells.Item(4,10).Value = "Active"
$ws1.Cells.Item(4,3).Interior.Color = $redColor

# Row 5: GBPUSD BUY
$ws1.Cells.Item(5,1).Value = "2025-07-28 19:11"
$ws1.Cells.Item(5,2).Value = "GBPUSD"
$ws1.Cells.Item(5,3).Value = "BUY"
$ws1.Cells.Item(5,4).Value = 1.26699
$ws1.Cells.Item(5,5).Value = 1.26497
$ws1.Cells.Item(5,6).Value = 1.27588
$ws1.Cells.Item(5,7).Value = 0.08
$ws1.Cells.Item(5,8).NumberFormat = "@"
$ws1.Cells.Item(5,8).Value = "71.0%"
$ws1.Cells.Item(5,9).Value = 4.42
$ws1.Cells.Item(5,10).Value = "Active"
$ws1.Cells.Item(5,3).Interior.Color = $greenColor

# Row 6: USDJPY SELL
$ws1.Cells.Item(6,1).Value = "2025-07-28 19:08"
$ws1.Cells.Item(6,2).Value = "USDJPY"
$ws1.Cells.Item(6,3).Value = "SELL"
$ws1.Cells.Item(6,4).Value = 149.11418
$ws1.Cells.Item(6,5).Value = 149.31518
$ws1.Cells.Item(6,6).Value = 148.20714
$ws1.Cells.Item(6,7).Value = 0.06
$ws1.Cells.Item(6,8).NumberFormat = "@"
$ws1.Cells.Item(6,8).Value = "66.0%"
$ws1.Cells.Item(6,9).Value = 4.51
$ws1.Cells.Item(6,10).Value = "Active"
$ws1.Cells.Item(6,3).Interior.Color = $redColor

# ===================== Sheet: Summary Dashboard =====================
$ws2 = $wb.Worksheets.Item("Summary Dashboard")

$ws2.Cells.Item(4,2).Value = 5
$ws2.Cells.Item(5,2).Value = 6
$ws2.Cells.Item(6,2).Value = 9
$ws2.Cells.Item(7,2).NumberFormat = "@"
$ws2.Cells.Item(7,2).Value = "79.2%"
$ws2.Cells.Item(8,2).NumberFormat = "@"
$ws2.Cells.Item(8,2).Value = "2.40"
$ws2.Cells.Item(9,2).NumberFormat = "@"
$ws2.Cells.Item(9,2).Value = "2025-07-28 19:35:29"

# ===================== Sheet: Signal History =====================
$ws3 = $wb.Worksheets.Item("Signal History")

# Row 2: AUDUSD SELL
$ws3.Cells.Item(2,1).Value = "2025-07-28 19:46"
$ws3.Cells.Item(2,2).Value = "AUDUSD"
$ws3.Cells.Item(2,3).Value = "SELL"
$ws3.Cells.Item(2,4).Value = 0.65914
$ws3.Cells.Item(2,5).Value = 0.66237
$ws3.Cells.Item(2,6).Value = 0.6515
$ws3.Cells.Item(2,7).Value = 0.1
$ws3.Cells.Item(2,8).Value = 0.74
$ws3.Cells.Item(2,9).Value = 2.37
$ws3.Cells.Item(2,10).Value = "Active"

# Row 3: GBPUSD SELL
$ws3.Cells.Item(3,1).Value = "2025-07-28 19:55"
$ws3.Cells.Item(3,2).Value = "GBPUSD"
$ws3.Cells.Item(3,3).Value = "SELL"
$ws3.Cells.Item(3,4).Value = 1.2672
$ws3.Cells.Item(3,5).Value = 1.27119
$ws3.Cells.Item(3,6).Value = 1.25817
$ws3.Cells.Item(3,7).Value = 0.09
$ws3.Cells.Item(3,8).Value = 0.86
$ws3.Cells.Item(3,9).Value = 2.27
$ws3.Cells.Item(3,10).Value = "Filled"

# Row 4: USDCAD SELL
$ws3.Cells.Item(4,1).Value = "2025-07-28 19:34"
$ws3.Cells.Item(4,2).Value = "USDCAD"
$ws3.Cells.Item(4,3).Value = "SELL"
$ws3.Cells.Item(4,4).Value = 1.36221
$ws3.Cells.Item(4,5).Value = 1.36531
$ws3.Cells.Item(4,6).Value = 1.3559
$ws3.Cells.Item(4,7).Value = 0.07
$ws3.Cells.Item(4,8).Value = 0.86
$ws3.Cells.Item(4,9).Value = 2.03
$ws3.Cells.Item(4,10).Value = "Active"

# Row 5: AUDUSD SELL
$ws3.Cells.Item(5,1).Value = "2025-07-28 19:45"
$ws3.Cells.Item(5,2).Value = "AUDUSD"
$ws3.Cells.Item(5,3).Value = "SELL"
$ws3.Cells.Item(5,4).Value = 0.65622
$ws3.Cells.Item(5,5).Value = 0.6601
$ws3.Cells.Item(5,6).Value = 0.64627
$ws3.Cells.Item(5,7).Value = 0.07
$ws3.Cells.Item(5,8).Value = 0.91
$ws3.Cells.Item(5,9).Value = 2.56
$ws3.Cells.Item(5,10).Value = "Pending"

# Row 6: NZDUSD SELL
$ws3.Cells.Item(6,1).Value = "2025-07-28 19:15"
$ws3.Cells.Item(6,2).Value = "NZDUSD"
$ws3.Cells.Item(6,3).Value = "SELL"
$ws3.Cells.Item(6,4).Value = 0.58906
$ws3.Cells.Item(6,5).Value = 0.59294
$ws3.Cells.Item(6,6).Value = 0.5839
$ws3.Cells.Item(6,7).Value = 0.01
$ws3.Cells.Item(6,8).Value = 0.77
$ws3.Cells.Item(6,9).Value = 1.33
$ws3.Cells.Item(6,10).Value = "Filled"

# Row 7: USDCHF BUY
$ws3.Cells.Item(7,1).Value = "2025-07-28 20:02"
$ws3.Cells.Item(7,2).Value = "USDCHF"
$ws3.Cells.Item(7,3).Value = "BUY"
$ws3.Cells.Item(7,4).Value = 0.88078
$ws3.Cells.Item(7,5).Value = 0.87585
$ws3.Cells.Item(7,6).Value = 0.88957
$ws3.Cells.Item(7,7).Value = 0.04
$ws3.Cells.Item(7,8).Value = 0.83
$ws3.Cells.Item(7,9).Value = 1.78
$ws3.Cells.Item(7,10).Value = "Filled"

# Row 8: USDCAD BUY
$ws3.Cells.Item(8,1).Value = "2025-07-28 19:59"
$ws3.Cells.Item(8,2).Value = "USDCAD"
$ws3.Cells.Item(8,3).Value = "BUY"
$ws3.Cells.Item(8,4).Value = 1.3664
$ws3.Cells.Item(8,5).Value = 1.36228
$ws3.Cells.Item(8,6).Value = 1.37606
$ws3.Cells.Item(8,7).Value = 0.08
$ws3.Cells.Item(8,8).Value = 0.75
$ws3.Cells.Item(8,9).Value = 2.34
$ws3.Cells.Item(8,10).Value = "Active"

# Row 9: NZDUSD SELL
$ws3.Cells.Item(9,1).Value = "2025-07-28 19:54"
$ws3.Cells.Item(9,2).Value = "NZDUSD"
$ws3.Cells.Item(9,3).Value = "SELL"
$ws3.Cells.Item(9,4).Value = 0.58913
$ws3.Cells.Item(9,5).Value = 0.59282
$ws3.Cells.Item(9,6).Value = 0.57921
$ws3.Cells.Item(9,7).Value = 0.08
$ws3.Cells.Item(9,8).Value = 0.89
$ws3.Cells.Item(9,9).Value = 2.69
$ws3.Cells.Item(9,10).Value = "Active"

# Row 10: NZDUSD BUY
$ws3.Cells.Item(10,1).Value = "2025-07-28 19:27"
$ws3.Cells.Item(10,2).Value = "NZDUSD"
$ws3.Cells.Item(10,3).Value = "BUY"
$ws3.Cells.Item(10,4).Value = 0.58736
$ws3.Cells.Item(10,5).Value = 0.58366
$ws3.Cells.Item(10,6).Value = 0.59299
$ws3.Cells.Item(10,7).Value = 0.05
$ws3.Cells.Item(10,8).Value = 0.8
$ws3.Cells.Item(10,9).Value = 1.52
$ws3.Cells.Item(10,10).Value = "Pending"

# Row 11: EURUSD SELL
$ws3.Cells.Item(11,1).Value = "2025-07-28 19:37"
$ws3.Cells.Item(11,2).Value = "EURUSD"
$ws3.Cells.Item(11,3).Value = "SELL"
$ws3.Cells.Item(11,4).Value = 1.09978
$ws3.Cells.Item(11,5).Value = 1.10314
$ws3.Cells.Item(11,6).Value = 1.09483
$ws3.Cells.Item(11,7).Value = 0.08
$ws3.Cells.Item(11,8).Value = 0.7
$ws3.Cells.Item(11,9).Value = 1.48
$ws3.Cells.Item(11,10).Value = "Pending"

# Row 12: AUDUSD SELL
$ws3.Cells.Item(12,1).Value = "2025-07-28 20:00"
$ws3.Cells.Item(12,2).Value = "AUDUSD"
$ws3.Cells.Item(12,3).Value = "SELL"
$ws3.Cells.Item(12,4).Value = 0.65556
$ws3.Cells.Item(12,5).Value = 0.65854
$ws3.Cells.Item(12,6).Value = 0.64765
$ws3.Cells.Item(12,7).Value = 0.02
$ws3.Cells.Item(12,8).Value = 0.77
$ws3.Cells.Item(12,9).Value = 2.66
$ws3.Cells.Item(12,10).Value = "Pending"

# Row 13: USDCAD BUY
$ws3.Cells.Item(13,1).Value = "2025-07-28 19:12"
$ws3.Cells.Item(13,2).Value = "USDCAD"
$ws3.Cells.Item(13,3).Value = "BUY"
$ws3.Cells.Item(13,4).Value = 1.36078
$ws3.Cells.Item(13,5).Value = 1.3558
$ws3.Cells.Item(13,6).Value = 1.36569
$ws3.Cells.Item(13,7).Value = 0.06
$ws3.Cells.Item(13,8).Value = 0.94
$ws3.Cells.Item(13,9).Value = 0.99
$ws3.Cells.Item(13,10).Value = "Pending"

# Row 14: GBPUSD BUY
$ws3.Cells.Item(14,1).Value = "2025-07-28 19:11"
$ws3.Cells.Item(14,2).Value = "GBPUSD"
$ws3.Cells.Item(14,3).Value = "BUY"
$ws3.Cells.Item(14,4).Value = 1.26699
$ws3.Cells.Item(14,5).Value = 1.26497
$ws3.Cells.Item(14,6).Value = 1.27588
$ws3.Cells.Item(14,7).Value = 0.08
$ws3.Cells.Item(14,8).Value = 0.71
$ws3.Cells.Item(14,9).Value = 4.42
$ws3.Cells.Item(14,10).Value = "Active"

# Row 15: USDJPY SELL
$ws3.Cells.Item(15,1).Value = "2025-07-28 19:08"
$ws3.Cells.Item(15,2).Value = "USDJPY"
$ws3.Cells.Item(15,3).Value = "SELL"
$ws3.Cells.Item(15,4).Value = 149.11418
$ws3.Cells.Item(15,5).Value = 149.31518
$ws3.Cells.Item(15,6).Value = 148.20714
$ws3.Cells.Item(15,7).Value = 0.06
$ws3.Cells.Item(15,8).Value = 0.66
$ws3.Cells.Item(15,9).Value = 4.51
$ws3.Cells.Item(15,10).Value = "Active"

# Row 16: USDCHF BUY
$ws3.Cells.Item(16,1).Value = "2025-07-28 19:52"
$ws3.Cells.Item(16,2).Value = "USDCHF"
$ws3.Cells.Item(16,3).Value = "BUY"
$ws3.Cells.Item(16,4).Value = 0.88032
$ws3.Cells.Item(16,5).Value = 0.87825
$ws3.Cells.Item(16,6).Value = 0.8867
$ws3.Cells.Item(16,7).Value = 0.07
$ws3.Cells.Item(16,8).Value = 0.69
$ws3.Cells.Item(16,9).Value = 3.08
$ws3.Cells.Item(16,10).Value = "Active"

Write-Output "done"
